$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 24, shifting existing rows 24-35 down to 25-36
$ws.Rows.Item(24).Insert()

$ws.Cells.Item(24, 1).Value = 5
$ws.Cells.Item(24, 2).Value = "Técnico `nÓtimo relacionamento com o cliente "
$ws.Cells.Item(24, 3).Value = 46043.57604188658
$ws.Cells.Item(24, 4).Value = "NjBmNTZjNjctMjFmYy00ZjU4LTg0NTItZDViMWVlZTE3NDUwOjU3MDE2"

# Reset the row height back to default (the multi-line comment otherwise
# leaves a custom auto-fit height on the new row that the original file
# doesn't have).
$ws.Rows.Item(24).EntireRow.AutoFit()
